$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, shifting existing rows 53-152 down to 54-153.
$ws.Rows(53).Insert()

# Populate the newly inserted row 53 with the new weekly price record.
$ws.Range("A53").Value = 10
$ws.Range("B53").Value = "Vega Modelo de Temuco"
$ws.Range("C53").Value = "La Araucanía"
$ws.Range("D53").Value = 44469
$ws.Range("E53").Value = 9
$ws.Range("F53").Value = 100112039
$ws.Range("G53").Value = "Ciboulette"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 30
$ws.Range("K53").Value = 7000
$ws.Range("L53").Value = 7000
$ws.Range("M53").Value = 7000
$ws.Range("N53").Value = "`$/docena de atados"
$ws.Range("O53").Value = "Provincia de Cautín"
$ws.Range("P53").Value = 2333
$ws.Range("Q53").Value = 3
$ws.Range("R53").Value = "Hortaliza"
